$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.789.86'
$ws.Range("E2").Value = '  -0.64%  '
$ws.Range("D3").Value = '1.598.16'
$ws.Range("E3").Value = '  -2.22%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = "'208.58"
$ws.Range("E5").Value = '  -2.58%  '
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("D7").Value = "'0.475"
$ws.Range("E7").Value = '  -5.88%  '
$ws.Range("E8").Value = '  -2.88%  '
$ws.Range("D9").Value = "'0.0607"
$ws.Range("E9").Value = '  -2.41%  '
$ws.Range("D10").Value = "'17.79"
$ws.Range("E10").Value = '  -3.64%  '
$ws.Range("D11").Value = "'0.0787"
$ws.Range("B12").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C12").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D12").Value = '1.820.20'
$ws.Range("E12").Value = '  -2.23%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.606.35'
$ws.Range("E13").Value = '  -1.92%  '
$ws.Range("D14").Value = "'4.03"
$ws.Range("D15").Value = "'0.508"
$ws.Range("E15").Value = '  -4.46%  '
$ws.Range("D16").Value = '25.786.16'
$ws.Range("E16").Value = '  -0.78%  '
$ws.Range("D17").Value = "'60.21"
$ws.Range("E17").Value = '  -2.28%  '
$ws.Range("E19").Value = '  +0.13%  '
$ws.Range("D20").Value = "'188.87"
$ws.Range("E20").Value = '  -0.78%  '
$ws.Range("E21").Value = '  -1.54%  '
$ws.Range("D22").Value = "'9.26"
$ws.Range("E22").Value = '  -3.33%  '
$ws.Range("E23").Value = '  -3.19%  '
$ws.Range("E24").Value = '  +0.17%  '
$ws.Range("D25").Value = "'141.13"
$ws.Range("E25").Value = '  -1.54%  '
$ws.Range("D26").Value = "'0.127"
$ws.Range("E26").Value = '  -3.87%  '
$ws.Range("D27").Value = "'1.71"
$ws.Range("E27").Value = '  -3.00%  '
$ws.Range("D28").Value = "'6.50"
$ws.Range("E28").Value = '  -4.27%  '
$ws.Range("D29").Value = "'14.84"
$ws.Range("E29").Value = '  -2.45%  '
$ws.Range("D30").Value = "'1.19"
$ws.Range("E30").Value = '  -3.42%  '
$ws.Range("E31").Value = '  -4.28%  '
$ws.Range("E32").Value = '  -2.72%  '
$ws.Range("D33").Value = "'2.98"
$ws.Range("E33").Value = '  -5.31%  '
$ws.Range("E34").Value = '  -1.32%  '
$ws.Range("E35").Value = '  -2.04%  '
$ws.Range("D36").Value = '1.092.30'
$ws.Range("E36").Value = '  -4.13%  '
$ws.Range("E37").Value = '  -3.18%  '
$ws.Range("D38").Value = "'0.791"
$ws.Range("E38").Value = '  -8.52%  '
$ws.Range("D39").Value = "'0.0150"
$ws.Range("E39").Value = '  -2.88%  '
$ws.Range("D40").Value = "'0.490"
$ws.Range("E40").Value = '  -6.44%  '
$ws.Range("D41").Value = "'95.55"
$ws.Range("E41").Value = '  -2.94%  '
$ws.Range("D42").Value = '1.733.55'
$ws.Range("D43").Value = "'5.08"
$ws.Range("E43").Value = '  -2.85%  '
$ws.Range("E44").Value = '  -5.36%  '
$ws.Range("D45").Value = "'52.96"
$ws.Range("E45").Value = '  -3.98%  '
$ws.Range("D46").Value = "'0.0510"
$ws.Range("E46").Value = '  -3.38%  '
$ws.Range("D47").Value = '0.0₇0973'
$ws.Range("E47").Value = '  -14.64%  '
$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").Value = "'0.411"
$ws.Range("E48").Value = '  -0.78%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").Value = "'1.42"
$ws.Range("E49").Value = '  -3.61%  '
$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = '  -0.20%  '
$ws.Range("D51").Value = "'7.35"
$ws.Range("E51").Value = '  -2.82%  '
